{"js": "// Insert a new \"Subtitle\" styled paragraph right after the document's\n// Title paragraph (\"Z\u00e1pis dok. jednotek\"), containing the descriptive\n// subtitle text \"Pokyny pro volbu a z\u00e1pis spr\u00e1vn\u00e9ho typu dokumenta\u010dn\u00edch\n// jednotek.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\n// Locate the Title paragraph (first paragraph of the document uses the\n// built-in \"Title\" style).\nlet titleParagraph = paragraphs.items[0];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].style === \"Title\") {\n    titleParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\n// Insert the new paragraph right after the title.\nconst subtitleParagraph = titleParagraph.insertParagraph(\n  \"Pokyny pro volbu a z\u00e1pis spr\u00e1vn\u00e9ho typu dokumenta\u010dn\u00edch jednotek.\",\n  Word.InsertLocation.after\n);\nsubtitleParagraph.style = \"Subtitle\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"Subtitle\" styled paragraph right after the document's\n# Title paragraph (\"Z\u00e1pis dok. jednotek\"), containing the descriptive\n# subtitle text \"Pokyny pro volbu a z\u00e1pis spr\u00e1vn\u00e9ho typu dokumenta\u010dn\u00edch\n# jednotek.\"\n\n$d = $word.ActiveDocument\n\n# Locate the Title paragraph (first paragraph of the document uses the\n# built-in \"Title\" style).\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n        break\n    }\n}\nif ($titlePara -eq $null) {\n    $titlePara = $d.Paragraphs(1)\n}\n\n# Split a new, empty paragraph in right after the title paragraph.\n$titlePara.Range.InsertParagraphAfter()\n\n# Grab that freshly created paragraph and fill it in.\n$subtitlePara = $titlePara.Next()\n$subtitlePara.Range.Text = \"Pokyny pro volbu a z\u00e1pis spr\u00e1vn\u00e9ho typu dokumenta\u010dn\u00edch jednotek.\"\n$subtitlePara.Style = \"Subtitle\"\n"}
